# Updates to mapping files. (Remove shipping and aviation, other changes.)
# Clears the "Air" and "Marine" transport-sector mapping rows (and the
# related note) on the "map" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("map")
$ws.Activate()

$ws.Range("B36").ClearContents()
$ws.Range("D36").ClearContents()
$ws.Range("B37").ClearContents()
$ws.Range("B44").ClearContents()
$ws.Range("C44").ClearContents()

$ws.Range("D37").Select()
